# Applies the "Commit izlazak s autoputa" edit:
#   1) "Vezani zahtjevi" -> "Vezani zahtj" | <bookmark _GoBack> | "evi"
#   2) "2. Izvrsavanje transakcije u odnosu na predjenu kilometrazu i tip vozila"
#        -> "2. Izvrsavanje transakcije u " | "odnosu na predjenu kilometrazu"
#        (the trailing " i tip vozila" is dropped)
#   3) the trailing empty paragraph that used to carry the _GoBack bookmark
#      loses that bookmark (it becomes a plain empty paragraph), because the
#      bookmark moved into the table cell edited in step 1.

$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# Helper pattern used below: Word (and this host) silently re-merges two
# abutting runs that end up with identical formatting, so simply writing
# text into a sub-range does not keep a persistent run boundary. Nudging the
# font color away and immediately back forces a genuine, lasting run split
# while leaving the run's formatting byte-for-byte the way it started.
# --------------------------------------------------------------------------

# --- Change 1: split "Vezani zahtjevi" after "Vezani zahtj" -------------
$rng1 = $d.Content
$rng1.Find.Execute("Vezani zahtjevi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos1 = $rng1.Start + 12

$nudge1 = $d.Range($splitPos1, $splitPos1 + 3)
$nudge1.Font.Color = 255
$restore1 = $d.Range($splitPos1, $splitPos1 + 3)
$restore1.Font.Color = 0

# --- Change 2: split the transaction sentence and drop " i tip vozila" ---
$rng2 = $d.Content
$rng2.Find.Execute("2. Izvršavanje transakcije u odnosu na pređenu kilometražu i tip vozila", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $rng2.Start
$splitPos2 = $start2 + 29

$nudge2 = $d.Range($splitPos2, $splitPos2 + 29)
$nudge2.Font.Color = 255
$restore2 = $d.Range($splitPos2, $splitPos2 + 29)
$restore2.Font.Color = 0

$tail2 = $d.Range($splitPos2 + 29, $start2 + 71)
$tail2.Text = ""

# --- Change 3: move the _GoBack bookmark from the trailing empty --------
# --- paragraph to the split point created in Change 1 -------------------
$d.Bookmarks("_GoBack").Delete()

$bmRange1 = $d.Range($splitPos1, $splitPos1)
$d.Bookmarks.Add("_GoBack", $bmRange1)

Write-Output "Edit applied."
